# Updates cryptos list values (Price / Volume(1h) columns) to match
# the latest scrape, per the GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '24.508.40'
$ws.Range('E2').Value = '  -1.00%  '
# Row 3
$ws.Range('D3').Value = '1.696.11'
$ws.Range('E3').Value = '  -0.32%  '
# Row 4
$ws.Range('E4').Value = '  -0.30%  '
# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '316.13'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.15%  '
# Row 6
$ws.Range('E6').Value = '  -0.33%  '
# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3922'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.40%  '
# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4075'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +0.80%  '
# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.492'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -2.46%  '
# Row 10
$ws.Range('E10').Value = '  -0.27%  '
# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '52.62'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -2.28%  '
# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08811'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -0.80%  '
# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '26.94'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +13.68%  '
# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.531'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +1.10%  '
# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '8.177'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -0.35%  '
# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.00001351'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +1.96%  '
# Row 17
$ws.Range('D17').Value = '1.690.20'
$ws.Range('E17').Value = '  -1.24%  '
# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '97.96'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -1.84%  '
# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.07173'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +1.57%  '
# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '20.71'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +5.03%  '
# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.307'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +3.12%  '
# Row 22
$ws.Range('E22').Value = '  -0.59%  '
# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '14.38'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -2.54%  '
# Row 24
$ws.Range('D24').Value = '24.497.79'
$ws.Range('E24').Value = '  -0.99%  '
# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.039'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -3.79%  '
# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.328'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -1.72%  '
# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.77'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +0.03%  '
# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '167.82'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +2.85%  '
# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.506'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -2.64%  '
# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '144.69'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +6.44%  '
# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.387'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +4.03%  '
# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.216'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +12.03%  '
# Row 33
$ws.Range('D33').Value = '1.875.89'
$ws.Range('E33').Value = '  -1.13%  '
# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08771'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -3.67%  '
# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '7.321'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -4.55%  '
# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.043'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -2.44%  '
# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.03057'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +9.86%  '
# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2808'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +1.76%  '
# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '10.94'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -1.74%  '
# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.09191'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +0.36%  '
# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '14.22'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -2.30%  '
# Row 42
$ws.Range('E42').Value = '  +4.24%  '
# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.479'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +1.13%  '
# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '17.47'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +9.51%  '
# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.664'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +3.12%  '
# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.7288'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +1.38%  '
# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.273'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +1.26%  '
# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.418'
$ws.Range('D48').ClearFormats()
# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '141.40'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +0.62%  '
# Row 51
$ws.Range('E51').Value = '  +2.06%  '
